$d = $word.ActiveDocument

# --- Step 1: insert the two new "MUG" paragraphs before the final (empty)
# paragraph, matching the Arial/shaded style used by the other database
# entries ("faceDB:" etc.) in this document. ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insPoint = $lastPara.Range
$insPoint.Collapse(1)

$newParasXml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="333333"/><w:szCs w:val="21"/><w:shd w:val="clear" w:color="auto" w:fill="F5F5F5"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="333333"/><w:szCs w:val="21"/><w:shd w:val="clear" w:color="auto" w:fill="F5F5F5"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial" w:hint="eastAsia"/><w:color w:val="333333"/><w:szCs w:val="21"/><w:shd w:val="clear" w:color="auto" w:fill="F5F5F5"/></w:rPr><w:t>M</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="333333"/><w:szCs w:val="21"/><w:shd w:val="clear" w:color="auto" w:fill="F5F5F5"/></w:rPr><w:t>UG Database:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$insPoint.InsertXML($newParasXml)

# --- Step 2: populate the (still) final paragraph with the MUG citation
# text, keeping its original paragraph formatting (Yu Mincho / ja-JP).
# Runs are appended one at a time so each becomes its own <w:r>, with
# w:proofErr spell-check markers bracketing the proper nouns, matching
# how Word would record it. ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $lastPara.Range
$r.Collapse(1)

$citationXml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:eastAsia="Yu Mincho" w:hint="eastAsia"/><w:lang w:eastAsia="ja-JP"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">N. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Aifanti</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, C. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Papachristou</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and A. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Delopoulos</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>,&#8221; The</w:t></w:r><w:r><w:t xml:space="preserve"> MUG Facial Expression Database,&#8221; in Proc. 11th Int. Workshop on Image Analysis for Multimedia Interactive Services (WIAMIS), </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Desenzano</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, Italy, April 12-14 2010.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r.InsertXML($citationXml)
